$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.963.66'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.375.35'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.94%  '

$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.676'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.94'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.71%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.47'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.74%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.40%  '

$ws.Range("E10").Value = '  +2.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.56'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '37.56'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +15.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.34'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.108'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.46'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.928'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.375.17'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.940.39'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("E19").Value = '  +2.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.22'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.61'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '254.95'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.22%  '

$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.77'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.97%  '

$ws.Range("E25").Value = '  -3.45%  '

$ws.Range("E26").Value = '  +0.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.60'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.21%  '

$ws.Range("E28").Value = '  +0.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.47'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '176.23'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.55%  '

$ws.Range("E31").Value = '  +1.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.135'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0759'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.43'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.66%  '

$ws.Range("E35").Value = '  -2.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.84'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.66'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.85%  '

$ws.Range("E38").Value = '  +3.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0282'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.43'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +15.19%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '20.80'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +10.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.19'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +12.62%  '

$ws.Range("E43").Value = '  -2.21%  '

$ws.Range("E44").Value = '  -1.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.12'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.58'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.79%  '

$ws.Range("E47").Value = '  +0.58%  '

$ws.Range("E48").Value = '  +0.06%  '

$ws.Range("E49").Value = '  -0.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '98.85'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.45'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +16.87%  '

Write-Host "Applied cryptos update"